$d = $word.ActiveDocument

# 1) Update the consignment date in the header paragraph.
$d.Content.Find.Execute("31/03/2022", $false, $false, $false, $false, $false,
                         $true, 1, $false, "09/05/2022", 2) | Out-Null

# 2) Fill in the order table (first/only table in the document).
$t = $d.Tables.Item(1)

# Row 2 ("#" = 1): quantity 1 -> 10, price 260 -> 250, sum 260 -> 2500.
$t.Cell(2, 3).Range.Text = "10"
$t.Cell(2, 5).Range.Text = "250"
$t.Cell(2, 6).Range.Text = "2500"

# Row 3 ("#" = 2): previously blank line item, now populated.
$t.Cell(3, 2).Range.Text = "Кефир 2,5% 900 гр."
$t.Cell(3, 3).Range.Text = "10"
$t.Cell(3, 4).Range.Text = "0"
$t.Cell(3, 5).Range.Text = "285"
$t.Cell(3, 6).Range.Text = "2850"
